$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.339.67"
$ws.Range("E2").Value = "  +1.31%  "

$ws.Range("D3").Value = "2.919.76"
$ws.Range("E3").Value = "  +4.27%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'353.57"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "'112.90"
$ws.Range("E6").Value = "  +3.60%  "

$ws.Range("D7").Value = "'0.561"
$ws.Range("E7").Value = "  +1.69%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.634"
$ws.Range("E9").Value = "  +0.86%  "

$ws.Range("D10").Value = "'40.31"
$ws.Range("E10").Value = "  +1.00%  "

$ws.Range("D11").Value = "'0.0866"
$ws.Range("E11").Value = "  +3.24%  "

$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("D13").Value = "'20.15"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("E14").Value = "  +1.58%  "

$ws.Range("D15").Value = "3.378.32"
$ws.Range("E15").Value = "  +4.33%  "

$ws.Range("E16").Value = "  +6.56%  "

$ws.Range("D17").Value = "2.926.71"
$ws.Range("E17").Value = "  +4.47%  "

$ws.Range("D18").Value = "52.368.51"
$ws.Range("E18").Value = "  +1.40%  "

$ws.Range("D19").Value = "'7.72"
$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").Value = "'3.33"
$ws.Range("E20").Value = "  +5.95%  "

$ws.Range("D21").Value = "'14.42"
$ws.Range("E21").Value = "  +6.66%  "

$ws.Range("E22").Value = "  +0.89%  "

$ws.Range("D23").Value = "'71.12"
$ws.Range("E23").Value = "  +1.09%  "

$ws.Range("D24").Value = "'271.47"
$ws.Range("E24").Value = "  +1.33%  "

$ws.Range("D25").Value = "'2.82"
$ws.Range("E25").Value = "  +2.82%  "

$ws.Range("D26").Value = "'26.92"
$ws.Range("E26").Value = "  +3.35%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("E28").Value = "  +1.28%  "

$ws.Range("D29").Value = "'10.65"
$ws.Range("E29").Value = "  +2.91%  "

$ws.Range("D30").Value = "'37.99"
$ws.Range("E30").Value = "  +1.47%  "

$ws.Range("D31").Value = "'6.53"
$ws.Range("E31").Value = "  +5.06%  "

$ws.Range("D32").Value = "'2.26"
$ws.Range("E32").Value = "  +1.38%  "

$ws.Range("D33").Value = "'6.17"
$ws.Range("E33").Value = "  +8.18%  "

$ws.Range("D34").Value = "'0.0963"
$ws.Range("E34").Value = "  +12.22%  "

$ws.Range("D35").Value = "'53.33"
$ws.Range("E35").Value = "  +2.70%  "

$ws.Range("D36").Value = "'0.0455"
$ws.Range("E36").Value = "  +2.24%  "

$ws.Range("D37").Value = "'0.999"

$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "  +6.18%  "

$ws.Range("D39").Value = "'18.89"
$ws.Range("E39").Value = "  +0.57%  "

$ws.Range("D40").Value = "'2.08"
$ws.Range("E40").Value = "  +4.19%  "

$ws.Range("E41").Value = "  +15.40%  "

$ws.Range("D42").Value = "'23.99"
$ws.Range("E42").Value = "  +9.59%  "

$ws.Range("D43").Value = "'0.118"
$ws.Range("E43").Value = "  +2.16%  "

$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'121.67"
$ws.Range("E44").Value = "  +2.29%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'2.62"
$ws.Range("E45").Value = "  +7.47%  "

$ws.Range("E46").Value = "  -0.43%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.56"
$ws.Range("E47").Value = "  +5.85%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.206.32"
$ws.Range("E48").Value = "  +4.91%  "

$ws.Range("E49").Value = "  +24.15%  "

$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").Value = "'0.972"
$ws.Range("E50").Value = "  +7.97%  "

$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "'0.0336"
$ws.Range("E51").Value = "  +12.93%  "

